# Apply May 2024 river trend-result refresh for Makakahi at end Kaiparoro Road
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0.992597607917238
$ws.Range("H2").Value = 0.822222222222222
$ws.Range("J2").Value = 2.15
$ws.Range("K2").Value = -0.360544240242516
$ws.Range("L2").Value = -0.569355314340893
$ws.Range("M2").Value = -0.161178957252842
$ws.Range("N2").Value = -16.7694995461635
$ws.Range("P2").Value = "Virtually certain improving"

# Row 3
$ws.Range("E3").Value = "WARNING: Sen slope based on tied non-censored values"
$ws.Range("F3").Value = 0.5
$ws.Range("H3").Value = 0.894736842105263
$ws.Range("J3").Value = 10.84
$ws.Range("K3").Value = 0.0
$ws.Range("L3").Value = -0.0935872714290756
$ws.Range("M3").Value = 0.0776525014842552
$ws.Range("N3").Value = 0.0
$ws.Range("P3").Value = "As likely as not increasing"

# Row 4
$ws.Range("F4").Value = 0.586853854667004
$ws.Range("G4").Value = 0.491525423728814
$ws.Range("M4").Value = 0.0
$ws.Range("P4").Value = "As likely as not improving"

# Row 5
$ws.Range("G5").Value = 0.925925925925926
$ws.Range("H5").Value = 0.111111111111111
$ws.Range("I5").Value = 2.0

# Row 6
$ws.Range("F6").Value = 0.956794793985594
$ws.Range("G6").Value = 0.728813559322034
$ws.Range("P6").Value = "Extremely likely improving"

# Row 7
$ws.Range("F7").Value = 0.0029835528725788
$ws.Range("H7").Value = 0.711864406779661
$ws.Range("J7").Value = 0.036
$ws.Range("K7").Value = 0.0054902178178963
$ws.Range("L7").Value = 0.002281008550135
$ws.Range("M7").Value = 0.0080819862828144
$ws.Range("N7").Value = 15.2506050497122
$ws.Range("P7").Value = "Exceptionally unlikely improving"

# Row 8
$ws.Range("F8").Value = 0.038387869874593
$ws.Range("H8").Value = 0.824561403508772
$ws.Range("J8").Value = 7.51
$ws.Range("K8").Value = -0.0543416275354601
$ws.Range("L8").Value = -0.104357142857143
$ws.Range("M8").Value = -0.0046658032749325
$ws.Range("N8").Value = -0.723590246810387
$ws.Range("P8").Value = "Extremely unlikely increasing"

# Row 9
$ws.Range("F9").Value = 0.0118824968841707
$ws.Range("H9").Value = 0.728813559322034
$ws.Range("J9").Value = 0.042
$ws.Range("K9").Value = 0.0061591738618524
$ws.Range("L9").Value = 0.0013664020889955
$ws.Range("M9").Value = 0.0090309065934065
$ws.Range("N9").Value = 14.6646996710772
$ws.Range("P9").Value = "Extremely unlikely improving"

# Row 10
$ws.Range("F10").Value = 0.0017618567376681
$ws.Range("H10").Value = 0.372881355932203
$ws.Range("J10").Value = 0.13
$ws.Range("K10").Value = 0.0153466386554622
$ws.Range("L10").Value = 0.0069642120843987
$ws.Range("M10").Value = 0.0266217201166181
$ws.Range("N10").Value = 11.8051066580478
$ws.Range("P10").Value = "Exceptionally unlikely improving"

# Row 11
$ws.Range("F11").Value = 0.0022781427562524
$ws.Range("H11").Value = 0.338983050847458
$ws.Range("J11").Value = 0.009
$ws.Range("K11").Value = 0.0008152901785714
$ws.Range("L11").Value = 0.000258778334363
$ws.Range("M11").Value = 0.0016042319481957
$ws.Range("N11").Value = 9.05877976190476
$ws.Range("P11").Value = "Exceptionally unlikely improving"

# Row 12
$ws.Range("F12").Value = 0.703021765672633
$ws.Range("H12").Value = 0.680412371134021
$ws.Range("J12").Value = 2.15
$ws.Range("K12").Value = -0.0358001543365546
$ws.Range("L12").Value = -0.154577793670439
$ws.Range("M12").Value = 0.0649056431029493
$ws.Range("N12").Value = -1.66512345751417
$ws.Range("P12").Value = "Likely improving"

# Row 13
$ws.Range("F13").Value = 0.312711884018881
$ws.Range("H13").Value = 0.846846846846847
$ws.Range("K13").Value = -0.0059325933946939
$ws.Range("L13").Value = -0.0372986628917775
$ws.Range("M13").Value = 0.0190568523946102
$ws.Range("N13").Value = -0.0547287213532654
$ws.Range("P13").Value = "Unlikely increasing"

# Row 14
$ws.Range("F14").Value = 0.999696112697561
$ws.Range("G14").Value = 0.368421052631579
$ws.Range("H14").Value = 0.0789473684210526
$ws.Range("L14").Value = -0.0003803102873802
$ws.Range("P14").Value = "Virtually certain improving"

# Row 15
$ws.Range("F15").Value = 0.433976145593665
$ws.Range("G15").Value = 0.943396226415094
$ws.Range("H15").Value = 0.113207547169811

# Row 16
$ws.Range("F16").Value = 0.703187014023289
$ws.Range("G16").Value = 0.640350877192982
$ws.Range("H16").Value = 0.0526315789473684
$ws.Range("P16").Value = "Likely improving"

# Row 17
$ws.Range("E17").Value = "ok"
$ws.Range("F17").Value = 0.000000032303912481976
$ws.Range("G17").Value = 0.087719298245614
$ws.Range("H17").Value = 0.543859649122807
$ws.Range("J17").Value = 0.027
$ws.Range("K17").Value = 0.0032783658787255
$ws.Range("L17").Value = 0.0022559094134296
$ws.Range("M17").Value = 0.0042654960227202
$ws.Range("N17").Value = 12.1420958471318

# Row 18
$ws.Range("F18").Value = 0.181302336457135
$ws.Range("H18").Value = 0.705357142857143
$ws.Range("J18").Value = 7.59
$ws.Range("K18").Value = -0.0108867603944563
$ws.Range("L18").Value = -0.0316249504025265
$ws.Range("M18").Value = 0.0087329949689335
$ws.Range("N18").Value = -0.143435578319583
$ws.Range("P18").Value = "Unlikely increasing"

# Row 19
$ws.Range("F19").Value = 0.00000150841796278942
$ws.Range("G19").Value = 0.0350877192982456
$ws.Range("H19").Value = 0.614035087719298
$ws.Range("J19").Value = 0.0351
$ws.Range("K19").Value = 0.0029873227917121
$ws.Range("L19").Value = 0.001994085532302
$ws.Range("M19").Value = 0.0039051162540533
$ws.Range("N19").Value = 8.51089114447893

# Row 20
$ws.Range("E20").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F20").Value = 0.0000000329746914724966
$ws.Range("G20").Value = 0.0265486725663717
$ws.Range("H20").Value = 0.309734513274336
$ws.Range("J20").Value = 0.11
$ws.Range("K20").Value = 0.0091266652083537
$ws.Range("L20").Value = 0.0062414559125085
$ws.Range("M20").Value = 0.0119995924568967
$ws.Range("N20").Value = 8.29696837123068

# Row 21
$ws.Range("F21").Value = 0.0007834317462165
$ws.Range("G21").Value = 0.0176991150442478
$ws.Range("H21").Value = 0.194690265486726
$ws.Range("J21").Value = 0.009
$ws.Range("K21").Value = 0.0003467014712861
$ws.Range("L21").Value = 0.0001317295015968
$ws.Range("M21").Value = 0.0005244077530509
$ws.Range("N21").Value = 3.85223856984654

# Rows 22-24: replace 5-year ASPM/MCI/QMCI figures with the 10-year figures
# that previously lived in rows 25-27 (trend period now matches C25:C27)
# Row 22
$ws.Range("C22").Value = 10.0
$ws.Range("F22").Value = 0.125726093254011
$ws.Range("J22").Value = 0.688
$ws.Range("K22").Value = -0.0108118257017981
$ws.Range("L22").Value = -0.0220206764173476
$ws.Range("M22").Value = 0.0019907599208247
$ws.Range("N22").Value = -1.57148629386601
$ws.Range("P22").Value = "Unlikely improving"

# Row 23
$ws.Range("C23").Value = 10.0
$ws.Range("F23").Value = 0.232756084937093
$ws.Range("J23").Value = 136.43
$ws.Range("K23").Value = -0.450194222344543
$ws.Range("L23").Value = -1.94705679077372
$ws.Range("M23").Value = 0.807008742925449
$ws.Range("N23").Value = -0.3299818385579

# Row 24
$ws.Range("C24").Value = 10.0
$ws.Range("F24").Value = 0.377227088747045
$ws.Range("K24").Value = -0.0096902357034891
$ws.Range("L24").Value = -0.110125628140703
$ws.Range("M24").Value = 0.0837510223125423
$ws.Range("N24").Value = -0.11903004180677

# Remove the now-duplicate old rows 25-27
$ws.Rows("25:27").Delete()

